$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.325.74'
$ws.Cells.Item(2, 5).Value = '  +0.66%  '
$ws.Cells.Item(3, 4).Value = '1.866.10'
$ws.Cells.Item(3, 5).Value = '  +0.45%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 5).Value = '  +0.11%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '235.87'
$ws.Cells.Item(5, 5).Value = '  +0.83%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.001'
$ws.Cells.Item(6, 5).Value = '  +0.12%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4675'
$ws.Cells.Item(7, 5).Value = '  -0.40%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2842'
$ws.Cells.Item(8, 5).Value = '  +0.83%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06531'
$ws.Cells.Item(9, 5).Value = '  -0.28%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '21.86'
$ws.Cells.Item(10, 5).Value = '  +8.54%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07932'
$ws.Cells.Item(11, 5).Value = '  +1.69%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '97.40'
$ws.Cells.Item(12, 5).Value = '  +0.14%  '
$ws.Cells.Item(13, 4).Value = '1.871.65'
$ws.Cells.Item(13, 5).Value = '  +0.77%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.158'
$ws.Cells.Item(14, 5).Value = '  +1.44%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6786'
$ws.Cells.Item(15, 5).Value = '  +1.12%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '279.84'
$ws.Cells.Item(16, 5).Value = '  -1.88%  '
$ws.Cells.Item(17, 4).Value = '30.332.06'
$ws.Cells.Item(17, 5).Value = '  +0.63%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '13.21'
$ws.Cells.Item(18, 5).Value = '  +4.72%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '1.001'
$ws.Cells.Item(19, 5).Value = '  +0.11%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '5.417'
$ws.Cells.Item(20, 5).Value = '  -0.39%  '
$ws.Cells.Item(21, 4).Value = '2.116.72'
$ws.Cells.Item(21, 5).Value = '  +0.87%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.000007315'
$ws.Cells.Item(22, 5).Value = '  +1.15%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '1.001'
$ws.Cells.Item(23, 5).Value = '  +0.11%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '6.154'
$ws.Cells.Item(24, 5).Value = '  +0.17%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '166.18'
$ws.Cells.Item(25, 5).Value = '  -1.02%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.166'
$ws.Cells.Item(26, 5).Value = '  -1.51%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '19.10'
$ws.Cells.Item(27, 5).Value = '  +0.21%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.937'
$ws.Cells.Item(28, 5).Value = '  +0.38%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.390'
$ws.Cells.Item(29, 5).Value = '  +3.70%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.09738'
$ws.Cells.Item(30, 5).Value = '  +0.96%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.406'
$ws.Cells.Item(31, 5).Value = '  +0.05%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.478'
$ws.Cells.Item(32, 5).Value = '  +0.58%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.090'
$ws.Cells.Item(33, 5).Value = '  -0.15%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.04732'
$ws.Cells.Item(34, 5).Value = '  +1.24%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.134'
$ws.Cells.Item(35, 5).Value = '  +4.14%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7091'
$ws.Cells.Item(36, 5).Value = '  +1.69%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.713'
$ws.Cells.Item(37, 5).Value = '  +0.43%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01869'
$ws.Cells.Item(38, 5).Value = '  +0.40%  '
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.576'
$ws.Cells.Item(39, 5).Value = '  +2.94%  '
$ws.Cells.Item(40, 2).Value = 'FraxShare'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.334'
$ws.Cells.Item(40, 5).Value = '  +0.28%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '74.80'
$ws.Cells.Item(41, 5).Value = '  +3.98%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.972'
$ws.Cells.Item(42, 5).Value = '  +1.80%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.8509'
$ws.Cells.Item(43, 5).Value = '  -1.22%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.4195'
$ws.Cells.Item(44, 5).Value = '  +0.74%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.001'
$ws.Cells.Item(45, 5).Value = '  +0.09%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '103.50'
$ws.Cells.Item(46, 5).Value = '  -0.71%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '972.14'
$ws.Cells.Item(47, 5).Value = '  -4.98%  '
$ws.Cells.Item(48, 5).Value = '  +3.41%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '7.206'
$ws.Cells.Item(49, 5).Value = '  -0.66%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '34.15'
$ws.Cells.Item(50, 5).Value = '  +1.09%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.1130'
$ws.Cells.Item(51, 5).Value = '  -1.15%  '
